$wb = $excel.ActiveWorkbook

# Add the new "ValidLogin" worksheet right after the existing "TC1" sheet
$tc1 = $wb.Worksheets.Item("TC1")
$newSheet = $wb.Worksheets.Add($null, $tc1)
$newSheet.Name = "ValidLogin"

# Re-fetch to be safe and populate the data-driven-testing table:
# header row (username/password) + a data row (admin/pointofsale)
$ws = $wb.Worksheets.Item("ValidLogin")
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"

# Make the new sheet the active / selected tab, matching the zoom level
# and selection state captured by the author
$ws.Select()
$excel.ActiveWindow.Zoom = 160
$ws.Range("B3").Select()
